$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price updates in column D (row number -> new value)
$dUpdates = @{
    2  = "283.59"
    4  = "6.206"
    5  = "0.06187"
    7  = "6.566"
    8  = "1.489"
    9  = "0.8179"
    10 = "0.01388"
    11 = "0.1650"
    12 = "0.08384"
    13 = "0.03665"
    14 = "0.03131"
    15 = "0.09130"
    16 = "3.698"
    17 = "0.001645"
    18 = "0.04676"
    19 = "0.006489"
    20 = "0.006194"
    22 = "0.0001499"
    23 = "3.794"
    25 = "0.3386"
    40 = "0.04722"
    45 = "0.00006407"
    47 = "0.9994"
    50 = "0.01239"
}

foreach ($row in $dUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$row]
    $cell.Style = "Normal"
}

# Rows 41-43 are being rotated / updated with new coin ordering and values.
# New row 41: CEJI
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "0.005597"
$d41.Style = "Normal"
$ws.Range("E41").Value = "40CEJICEJI"

# New row 42: KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "0.007073"
$d42.Style = "Normal"
$ws.Range("E42").Value = "41KickTokenKICK"

# New row 43: BKEXToken
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$d43 = $ws.Range("D43")
$d43.NumberFormat = "@"
$d43.Value = "0.1102"
$d43.Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
